$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18; this shifts current rows 18-38 down to 19-39
$ws.Rows.Item(18).Insert()

# Fill in the new row 18 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,J,K,R are constant across the whole sheet.
$ws.Cells.Item(18, 1).Value = 3
$ws.Cells.Item(18, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(18, 3).Value = "Coquimbo"
$ws.Cells.Item(18, 4).Value = 44469
$ws.Cells.Item(18, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(18, 5).Value = 5
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100108
$ws.Cells.Item(18, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(18, 9).Value = 100108004
$ws.Cells.Item(18, 10).Value = "Papaya"
$ws.Cells.Item(18, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 50
$ws.Cells.Item(18, 14).Value = 16000
$ws.Cells.Item(18, 15).Value = 16000
$ws.Cells.Item(18, 16).Value = 16000
$ws.Cells.Item(18, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(18, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(18, 19).Value = 1600
$ws.Cells.Item(18, 20).Value = 10
